$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product backlog")

# --- Row 6: update requirement 5 texts ---
$ws.Range("B6").Value = "Als systeem wil ik sensordata waterpas hebben, zodat ik altijd nauwkeurige en realistische lezingen krijg."
$ws.Range("C6").Value = "1. Onderzoek naar data-analyse en filtering`n2. Algoritme om sensor data waterpas te houden`n3. Algoritme testen en kalibreren"
$ws.Range("D6").Value = "Het algortime geeft accurate data die waterpas is en gehouden word."
$ws.Range("A6").Value = "5. Algoritme om constante data over de positie van de Floating Farm te krijgen"

# --- Row 7: clear leftover content + formatting in A:D ---
$ws.Range("A7:D7").Clear()

# --- Sheet view: scroll + selection change ---
$ws.Application.ActiveWindow.ScrollRow = 3
$ws.Range("C10").Select()
